$p = $ppt.ActivePresentation
$p.HasNotesMaster = $false
$p.HasNotesMaster = $true
$nm = $p.NotesMaster
Write-Output $nm.Name
